# Update column E (Obrigatorio) from "N" to "S" for the rows corresponding
# to the fields that became mandatory in this layout revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(2, 3, 4, 5, 6, 7, 8, 12, 13, 14, 17, 18, 19, 22)

foreach ($r in $rows) {
    $ws.Range("E$r").Value = "S"
}
